$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (current GORONGOSA row), shifting it down to row 3
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the CHIMANIMANI data
$ws.Range("A2").Value = "CHIMANIMANI"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 148
$ws.Range("D2").Value = 119
$ws.Range("E2").Value = 110
$ws.Range("F2").Value = 66
$ws.Range("G2").Value = 443
